# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.393.32"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "3.495.52"
$ws.Range("E3").Value = "  +6.01%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").Value = "'245.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.57%  "

$ws.Range("D6").Value = "'651.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.28%  "

$ws.Range("D7").Value = "'1.43"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.92%  "

$ws.Range("D8").Value = "'0.416"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.45%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("D10").Value = "'1.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.48%  "

$ws.Range("D11").Value = "3.487.38"
$ws.Range("E11").Value = "  +5.84%  "

$ws.Range("D12").Value = "'43.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.98%  "

$ws.Range("D13").Value = "'0.200"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("D14").Value = "'6.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.97%  "

$ws.Range("D15").Value = "96.144.27"
$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("D16").Value = "4.140.54"
$ws.Range("E16").Value = "  +5.85%  "

$ws.Range("D17").Value = "'0.0000255"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.67%  "

$ws.Range("D18").Value = "'8.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.69%  "

$ws.Range("D19").Value = "3.502.11"
$ws.Range("E19").Value = "  +6.21%  "

$ws.Range("D20").Value = "'18.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +13.38%  "

$ws.Range("D21").Value = "'12.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +20.16%  "

$ws.Range("D22").Value = "'0.500"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.70%  "

$ws.Range("D23").Value = "'516.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.81%  "

$ws.Range("D24").Value = "'3.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.61%  "

$ws.Range("D25").Value = "'0.0000195"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.82%  "

$ws.Range("D26").Value = "'6.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.67%  "

$ws.Range("D27").Value = "'92.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.81%  "

$ws.Range("D28").Value = "'12.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.98%  "

$ws.Range("D29").Value = "3.663.63"
$ws.Range("E29").Value = "  +5.15%  "

$ws.Range("D30").Value = "'12.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +15.48%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +17.45%  "

$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").Value = "'1.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("D33").Value = "'0.141"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.03%  "

$ws.Range("D34").Value = "'0.186"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.21%  "

$ws.Range("D35").Value = "'31.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +15.40%  "

$ws.Range("D36").Value = "'0.587"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.51%  "

$ws.Range("E37").Value = "  +0.45%  "

$ws.Range("D38").Value = "'7.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.99%  "

$ws.Range("D39").Value = "'1.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.77%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'0.943"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +17.90%  "

$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'521.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.42%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.153"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.23%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("E44").Value = "  -0.93%  "

$ws.Range("D45").Value = "'1.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.14%  "

$ws.Range("E46").Value = "  +8.39%  "

$ws.Range("D47").Value = "'3.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("D48").Value = "'5.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.06%  "

$ws.Range("D49").Value = "'3.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.66%  "

$ws.Range("D50").Value = "'2.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.04%  "

$ws.Range("D51").Value = "'8.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.18%  "
